$d = $word.ActiveDocument
# Scope the edit to the first paragraph only (it holds both the floating
# textbox and the title line) so the rest of the document - most notably the
# table - is left completely untouched by the XML round trip.
$rng = $d.Paragraphs(1).Range
$xml = $rng.WordOpenXML

# ---------------------------------------------------------------------------
# 1) Textbox caption: merge the two runs "一份交" + "给运送人便于核对" into a
#    single run and drop the now-empty _GoBack bookmark that used to sit
#    between them.
# ---------------------------------------------------------------------------
$oldCaption = '<w:r w:rsidRPr="000631C6"><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/></w:rPr><w:t>一份交</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidRPr="000631C6"><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/></w:rPr><w:t>给运送人便于核对</w:t></w:r>'
$newCaption = '<w:r w:rsidRPr="000631C6"><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/></w:rPr><w:t>一份交给运送人便于核对</w:t></w:r>'

if ($xml.Contains($oldCaption)) {
    $xml = $xml.Replace($oldCaption, $newCaption)
}

# ---------------------------------------------------------------------------
# 2) Title run "外协加工清单" becomes "外协【加工】清单", split across five
#    runs, with the _GoBack bookmark relocated between "】" and "清单".
# ---------------------------------------------------------------------------
$titleRpr = '<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:hint="eastAsia"/><w:sz w:val="28"/></w:rPr>'
$oldTitle = '<w:r>' + $titleRpr + '<w:t>外协加工清单</w:t></w:r>'
$newTitle = '<w:r>' + $titleRpr + '<w:t>外协</w:t></w:r>' + `
            '<w:r>' + $titleRpr + '<w:t>【</w:t></w:r>' + `
            '<w:r>' + $titleRpr + '<w:t>加工</w:t></w:r>' + `
            '<w:r>' + $titleRpr + '<w:t>】</w:t></w:r>' + `
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
            '<w:r>' + $titleRpr + '<w:t>清单</w:t></w:r>'

if ($xml.Contains($oldTitle)) {
    $xml = $xml.Replace($oldTitle, $newTitle)
}

# ---------------------------------------------------------------------------
# 3) "[" + "CreateDate" (w:proofErr-wrapped) + "]" collapse into one run's
#    text "[CreateDate]" (the engine already normalizes consecutive runs
#    with identical formatting and drops proofErr markers on XML round trip,
#    but we fix it up explicitly too in case that text still is split).
# ---------------------------------------------------------------------------
$dateRpr = '<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:hint="eastAsia"/><w:sz w:val="28"/></w:rPr>'
$oldDate = '<w:r w:rsidRPr="00B277D8">' + $dateRpr + '<w:t>[</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00B277D8">' + $dateRpr + '<w:t>CreateDate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00B277D8">' + $dateRpr + '<w:t>]</w:t></w:r>'
$newDate = '<w:r w:rsidRPr="00B277D8">' + $dateRpr + '<w:t>[CreateDate]</w:t></w:r>'

if ($xml.Contains($oldDate)) {
    $xml = $xml.Replace($oldDate, $newDate)
}

$rng.InsertXML($xml)
